# Update cryptos list: refresh Price (column D) and Volume(1h) (column E)
# values for each coin row, as produced by the scheduled GitHub Actions job.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "67.206.03"; E = "  -1.97%  " }
    @{ Row = 3; D = "2.669.00"; E = "  -1.05%  " }
    @{ Row = 4; D = "1.00"; E = "  -0.04%  " }
    @{ Row = 5; D = "596.11"; E = "  -0.49%  " }
    @{ Row = 6; D = "164.45"; E = "  +2.69%  " }
    @{ Row = 7; D = $null; E = "  +0.03%  " }
    @{ Row = 8; D = "0.544"; E = "  -0.08%  " }
    @{ Row = 9; D = "2.667.04"; E = "  -1.10%  " }
    @{ Row = 10; D = $null; E = "  +0.69%  " }
    @{ Row = 11; D = $null; E = "  +0.89%  " }
    @{ Row = 12; D = "0.356"; E = "  -1.24%  " }
    @{ Row = 13; D = "5.19"; E = "  -2.07%  " }
    @{ Row = 14; D = "27.61"; E = "  -2.35%  " }
    @{ Row = 15; D = "3.153.43"; E = "  -1.27%  " }
    @{ Row = 16; D = $null; E = "  -3.21%  " }
    @{ Row = 17; D = "67.147.66"; E = "  -1.95%  " }
    @{ Row = 18; D = "2.654.50"; E = "  -1.90%  " }
    @{ Row = 19; D = "11.65"; E = "  -2.72%  " }
    @{ Row = 20; D = "361.39"; E = "  -1.44%  " }
    @{ Row = 21; D = "7.50"; E = "  -2.48%  " }
    @{ Row = 22; D = $null; E = "  -4.18%  " }
    @{ Row = 23; D = "4.79"; E = "  -2.23%  " }
    @{ Row = 24; D = $null; E = "  -5.12%  " }
    @{ Row = 25; D = $null; E = "  +0.07%  " }
    @{ Row = 26; D = "71.01"; E = "  -4.71%  " }
    @{ Row = 27; D = "10.08"; E = "  -0.17%  " }
    @{ Row = 28; D = $null; E = "  -0.58%  " }
    @{ Row = 29; D = "0.998"; E = "  -0.12%  " }
    @{ Row = 30; D = $null; E = "  -3.35%  " }
    @{ Row = 31; D = "550.96"; E = "  -3.76%  " }
    @{ Row = 32; D = "7.95"; E = "  -3.44%  " }
    @{ Row = 33; D = $null; E = "  -5.70%  " }
    @{ Row = 34; D = "1.92"; E = "  -1.55%  " }
    @{ Row = 35; D = $null; E = "  -3.07%  " }
    @{ Row = 36; D = $null; E = "  -0.01%  " }
    @{ Row = 37; D = $null; E = "  -5.24%  " }
    @{ Row = 38; D = "19.45"; E = "  -2.62%  " }
    @{ Row = 39; D = "155.70"; E = "  -3.26%  " }
    @{ Row = 40; D = "0.372"; E = "  -2.23%  " }
    @{ Row = 41; D = $null; E = "  -2.89%  " }
    @{ Row = 42; D = "1.82"; E = "  -4.89%  " }
    @{ Row = 43; D = $null; E = "  +0.19%  " }
    @{ Row = 44; D = $null; E = "  +0.00%  " }
    @{ Row = 45; D = $null; E = "  -5.85%  " }
    @{ Row = 46; D = "40.22"; E = "  -0.61%  " }
    @{ Row = 47; D = "0.0₆0296"; E = "  -6.42%  " }
    @{ Row = 48; D = "0.585"; E = "  -2.63%  " }
    @{ Row = 49; D = "152.41"; E = "  -3.68%  " }
    @{ Row = 50; D = "3.82"; E = "  -4.08%  " }
    @{ Row = 51; D = $null; E = "  -3.77%  " }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($null -ne $u.D) {
        $priceCell = $ws.Cells.Item($row, 4)
        # Force the cell to remain plain text so price strings such as
        # "1.00" or "67.206.03" are not reinterpreted as numbers.
        $priceCell.NumberFormat = "@"
        $priceCell.Value = $u.D
    }

    if ($null -ne $u.E) {
        $ws.Cells.Item($row, 5).Value = $u.E
    }
}
